$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")

# Make sure the "Commands" sheet is the active tab (it was already the
# active tab in the source workbook).
$ws.Activate()

# --- Add the new "WebSocket Command" action (row 86) ---
$ws.Range("A86").Value = "WebSocket Command"
$ws.Range("B86").Value = "send(<json>)"
$ws.Range("C86").Value = "If {} substitutions are used, json brackets need to be duplicated to escape them like in send({{ " + [char]8220 + "value" + [char]8221 + ": {}}})"

# --- Add the "Probat Sample Roaster" sleep command row (row 87) ---
$ws.Range("B87").Value = "sleep(<float>)"
$ws.Range("C87").Value = "sleep: add a delay of <float> seconds"

# Update the selection to reflect the newly added rows, matching the
# author's edit (rows 86:87 selected, active cell A86).
$ws.Range("A86:C87").Select()
